# Updated cryptos list on Fri Jul  7 19:52:00 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (prices like
    # "1.0000" or "30.215.10") are not coerced into numbers/dates,
    # then restore General/Normal so no stray cell style lingers.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '30.215.10'
Set-TextValue 'E2' '  -0.27%  '
Set-TextValue 'D3' '1.863.86'
Set-TextValue 'E3' '  -1.24%  '
Set-TextValue 'D4' '0.9999'
Set-TextValue 'E4' '  -0.07%  '
Set-TextValue 'D5' '235.08'
Set-TextValue 'E5' '  -1.11%  '
Set-TextValue 'E6' '  -0.01%  '
Set-TextValue 'D7' '0.4662'
Set-TextValue 'E7' '  -0.43%  '
Set-TextValue 'D8' '0.2829'
Set-TextValue 'E8' '  -0.43%  '
Set-TextValue 'D9' '0.06524'
Set-TextValue 'E9' '  -1.04%  '
Set-TextValue 'D10' '21.41'
Set-TextValue 'E10' '  +3.01%  '
Set-TextValue 'D11' '0.07858'
Set-TextValue 'E11' '  +1.06%  '
Set-TextValue 'D12' '97.39'
Set-TextValue 'E12' '  -0.31%  '
Set-TextValue 'D13' '1.869.03'
Set-TextValue 'E13' '  -0.95%  '
Set-TextValue 'D14' '5.098'
Set-TextValue 'E14' '  +0.02%  '
Set-TextValue 'D15' '0.6722'
Set-TextValue 'E15' '  -0.54%  '
Set-TextValue 'D16' '280.23'
Set-TextValue 'E16' '  -1.37%  '
Set-TextValue 'D17' '30.192.89'
Set-TextValue 'E17' '  -0.41%  '
Set-TextValue 'D18' '0.9996'
Set-TextValue 'E18' '  -0.05%  '
Set-TextValue 'D19' '5.520'
Set-TextValue 'E19' '  +2.23%  '
Set-TextValue 'D20' '12.66'
Set-TextValue 'E20' '  +0.07%  '
Set-TextValue 'D21' '2.111.80'
Set-TextValue 'E21' '  -0.70%  '
Set-TextValue 'D22' '0.000007274'
Set-TextValue 'E22' '  -0.25%  '
Set-TextValue 'D23' '1.001'
Set-TextValue 'E23' '  +0.00%  '
Set-TextValue 'E24' '  -0.57%  '
Set-TextValue 'D25' '9.188'
Set-TextValue 'E25' '  -2.41%  '
Set-TextValue 'D26' '164.33'
Set-TextValue 'E26' '  -2.03%  '
Set-TextValue 'D27' '19.12'
Set-TextValue 'E27' '  -0.69%  '
Set-TextValue 'E28' '  -3.13%  '
Set-TextValue 'D29' '1.383'
Set-TextValue 'E29' '  +0.34%  '
Set-TextValue 'D30' '0.09690'
Set-TextValue 'E30' '  -0.70%  '
Set-TextValue 'D31' '4.412'
Set-TextValue 'E31' '  +0.93%  '
Set-TextValue 'E32' '  -0.57%  '
Set-TextValue 'E33' '  -1.33%  '
Set-TextValue 'D34' '0.04690'
Set-TextValue 'E34' '  +0.39%  '
Set-TextValue 'D35' '1.113'
Set-TextValue 'E35' '  +1.34%  '
Set-TextValue 'D36' '0.7061'
Set-TextValue 'E36' '  -0.18%  '
Set-TextValue 'D37' '2.728'
Set-TextValue 'E37' '  +0.55%  '
Set-TextValue 'D38' '0.01851'
Set-TextValue 'E38' '  -1.01%  '
Set-TextValue 'D39' '2.528'
Set-TextValue 'D40' '6.219'
Set-TextValue 'E40' '  -7.29%  '
Set-TextValue 'D41' '73.15'
Set-TextValue 'D42' '1.934'
Set-TextValue 'E42' '  -1.59%  '
Set-TextValue 'D43' '0.8462'
Set-TextValue 'E43' '  -2.62%  '
Set-TextValue 'B44' 'Quant'
Set-TextValue 'C44' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue 'D44' '103.98'
Set-TextValue 'E44' '  +0.23%  '
Set-TextValue 'B45' 'PaxDollar'
Set-TextValue 'C45' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D45' '1.0000'
Set-TextValue 'E45' '  +0.01%  '
Set-TextValue 'B46' 'TheSandbox'
Set-TextValue 'C46' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D46' '0.4160'
Set-TextValue 'E46' '  -0.54%  '
Set-TextValue 'D47' '7.186'
Set-TextValue 'E47' '  -1.27%  '
Set-TextValue 'D48' '9.177'
Set-TextValue 'E48' '  -0.32%  '
Set-TextValue 'D49' '934.34'
Set-TextValue 'E49' '  -5.80%  '
Set-TextValue 'D50' '34.12'
Set-TextValue 'E50' '  +0.48%  '
Set-TextValue 'E51' '  -2.06%  '
